# Update LR-pairs data with new TPM-derived values (columns G:T, rows 2-10)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(10.37574666666667,31.12724,0.8643482197679554,0.8643482197679554,2,0.6666666666666666,0.08106466666666667,0.243194,0.01252465659474717,0.01252465659474717,0.8411064449511112,7.56995800456,0.0108256646308747,0.0108256646308747),
    @(10.37574666666667,31.12724,0.8643482197679554,0.8643482197679554,3,1,5.993807666666666,17.981423,0.926055528343168,0.926055528343168,62.19022991805777,559.7120692625199,0.8004344473296906,0.8004344473296906),
    @(10.37574666666667,31.12724,0.8643482197679554,0.8643482197679554,2,0.6666666666666666,0.397534,1.192602,0.06141981506208485,0.06141981506208484,4.124712075386666,37.12240867847999,0.05308810780739009,0.05308810780739008),
    @(1.154570666666667,3.463712,0.09618113591146868,0.09618113591146869,2,0.6666666666666666,0.08106466666666667,0.243194,0.01252465659474717,0.01252465659474717,0.09359488623644445,0.8423539761280001,0.00120463569818385,0.00120463569818385),
    @(1.154570666666667,3.463712,0.09618113591146868,0.09618113591146869,3,1,5.993807666666666,17.981423,0.926055528343168,0.926055528343168,6.92027451357511,62.282470622176,0.08906907263314118,0.08906907263314119),
    @(1.154570666666667,3.463712,0.09618113591146868,0.09618113591146869,2,0.6666666666666666,0.397534,1.192602,0.06141981506208485,0.06141981506208484,0.4589810954026666,4.130829858624,0.005907427580143653,0.005907427580143654),
    @(0.4738106666666667,1.421432,0.03947064432057595,0.03947064432057595,2,0.6666666666666666,0.08106466666666667,0.243194,0.01252465659474717,0.01252465659474717,0.03840930375644445,0.345683733808,0.0004943562656886214,0.0004943562656886214),
    @(0.4738106666666667,1.421432,0.03947064432057595,0.03947064432057595,3,1,5.993807666666666,17.981423,0.926055528343168,0.926055528343168,2.839930006415111,25.559370057736,0.03655200838033622,0.03655200838033622),
    @(0.4738106666666667,1.421432,0.03947064432057595,0.03947064432057595,2,0.6666666666666666,0.397534,1.192602,0.06141981506208485,0.06141981506208484,0.1883558495626667,1.695202646064,0.002424279674551104,0.002424279674551104)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = 7 + $j   # column G = 7
        $ws.Cells.Item($r, $col).Value = $rowValues[$j]
    }
}
